$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.709.06'
$ws.Range("E2").Value = '  -2.27%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.875.32'
$ws.Range("E3").Value = '  -2.00%  '
$ws.Range("E4").Value = '  -0.83%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.98'
$ws.Range("E5").Value = '  -0.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.692'
$ws.Range("E6").Value = '  -4.24%  '
$ws.Range("E7").Value = '  -0.88%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.21'
$ws.Range("E8").Value = '  +1.37%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.347'
$ws.Range("E9").Value = '  -2.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '50.54'
$ws.Range("E10").Value = '  -4.82%  '
$ws.Range("E11").Value = '  +0.54%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0967'
$ws.Range("E12").Value = '  -2.14%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '12.92'
$ws.Range("E13").Value = '  +2.39%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.148.11'
$ws.Range("E14").Value = '  -1.92%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.715'
$ws.Range("E15").Value = '  -0.49%  '
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.895.32'
$ws.Range("E17").Value = '  -1.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '34.708.21'
$ws.Range("E18").Value = '  -2.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.75'
$ws.Range("E19").Value = '  -0.51%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0820'
$ws.Range("E20").Value = '  -0.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '245.44'
$ws.Range("E21").Value = '  +1.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.72'
$ws.Range("E22").Value = '  -3.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.90'
$ws.Range("E23").Value = '  -3.13%  '
$ws.Range("E24").Value = '  -0.85%  '
$ws.Range("E25").Value = '  +3.38%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.25'
$ws.Range("E26").Value = '  -4.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.75'
$ws.Range("E27").Value = '  -1.76%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.38'
$ws.Range("E28").Value = '  -3.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.21'
$ws.Range("E29").Value = '  -3.35%  '
$ws.Range("E30").Value = '  -4.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.128.40'
$ws.Range("E31").Value = '  -0.35%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.71'
$ws.Range("E32").Value = '  +16.53%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.22'
$ws.Range("E33").Value = '  -3.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0578'
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("E35").Value = '  -2.08%  '
$ws.Range("E36").Value = '  -0.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.71'
$ws.Range("E37").Value = '  -13.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.836'
$ws.Range("E38").Value = '  -8.84%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.98'
$ws.Range("E39").Value = '  -3.80%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.55'
$ws.Range("E40").Value = '  -2.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '98.01'
$ws.Range("E41").Value = '  -0.90%  '
$ws.Range("E42").Value = '  +0.17%  '
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("E44").Value = '  -4.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.289.24'
$ws.Range("E45").Value = '  -4.28%  '
$ws.Range("E46").Value = '  -5.96%  '
$ws.Range("E47").Value = '  -0.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.74'
$ws.Range("E48").Value = '  -1.45%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0767'
$ws.Range("E49").Value = '  +5.98%  '
$ws.Range("B50").Value = 'Gas'
$ws.Range("C50").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '12.10'
$ws.Range("E50").Value = '  +0.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.49'
$ws.Range("E51").Value = '  -1.20%  '
